$data = @(
    @(-0.8011660996346401, 0.5787938608431268, 0.1949980559329599, 0.8798839467072684, 0.1878585623255973, 0.4273069063282261, 0.4640348830873968, 1.290528957474712, 0.6347307473531471, -0.02341949237247376),
    @(0.9961641555676001, 1.681050046341908, 0.9890246619602374, 1.228473005962866, 1.265200982722037, 2.091695057109352, 1.435896846987787, 0.7777466072621664, 1.176861555113949, 1.115425012992372),
    @(-0.007139493607362657, 0.2323088503952662, 0.2690368271544369, 1.095530901541752, 0.4397326914201872, -0.2184175483054337, 0.1806973995463494, 0.119260857424772, -0.3149050776923316, 0.6638963143867453),
    @(0.2761763207617995, 1.102670395149115, 0.4468721850275498, -0.211278054698071, 0.187836893153712, 0.1264003510321347, -0.3077655840849689, 0.671035807994108, 0.4439199035461818, -0.1895398258799067),
    @(0.1706958642657503, -0.4874543754598706, -0.0883394276080875, -0.1497759697296649, -0.5839419048467684, 0.3948594872323085, 0.1677435827843823, -0.4657161466417062, -0.0383191879861019, -0.2046561426742488),
    @(-0.2590352918738378, -0.3204718339954152, -0.7546377691125187, 0.2241636229665582, -0.00295228148136796, -0.6364120109074565, -0.2090150522518522, -0.3753520069399991, -0.2526409904991733, -0.205407253944033),
    @(-0.4956024772386809, 0.483198914840396, 0.2560830103924698, -0.3773767190336187, 0.0500202396219856, -0.1163167150661613, 0.006394301374664513, 0.05362803792980481, -0.3429906935926468, -0.08893865846892673),
    @(0.7516854876311507, 0.1182257582050622, 0.5456227168606665, 0.3792857621725196, 0.5019967786133455, 0.5492305151684858, 0.1526117836460341, 0.4066638187697542, 0.6267882086433268, 0.3909997600566587),
    @(-0.2060627707704842, -0.3723997254586311, -0.2496887090178053, -0.202454972462665, -0.5990737039851166, -0.3450216688613965, -0.1248972789878239, -0.3606857275744921, -0.5396584850452768, -0.5502190153045978),
    @(-0.04362593824732108, 0.00360779830781921, -0.3930109332146324, -0.1389588980909123, 0.0811654917826603, -0.1546229568040079, -0.3335957142747926, -0.3441562445341136, -0.5518753546922437, -0.334666455217339),
    @(-0.3493849949673113, -0.09533295984359125, 0.1247914300299814, -0.1109970185566868, -0.2899697760274715, -0.3005303062867926, -0.5082494164449226, -0.2910405169700179, 0.0186456522880436, -0.1356584815439604),
    @(0.4741764249972927, 0.2383879764106245, 0.05941521893983981, 0.04885468868051879, -0.1588644214776113, 0.05834447799729348, 0.3680306472553549, 0.2137265134233509, 0.8641529346425341, 0.647987886024325),
    @(-0.4147612060574529, -0.4253217363167739, -0.633040846474904, -0.4158319469999993, -0.1061457777419378, -0.2604499115739418, 0.3899765096452414, 0.1738114610270322, -0.2332176680079241, 0.2658245491694957),
    @(-0.2182796404174511, -0.001070740942546333, 0.3086154283155151, 0.1543112944835111, 0.8047377157026943, 0.5885726670844852, 0.1815435380495288, 0.6805857552269486, 0.5217368420714361, 0.5779006630416801),
    @(0.5268950687329662, 0.3725909349009622, 1.023017356120145, 0.8068523075019363, 0.3998231784669799, 0.8988653956443997, 0.7400164824888872, 0.7961803034591312, 3.045983143070524, 10.51513347392782),
    @(0.1553820354260574, 0.8058084566452406, 0.5896434080270315, 0.1826142789920752, 0.681656496169495, 0.5228075830139824, 0.5789714039842264, 2.828774243595618, 10.29792457445291, -7.873337225684557),
    @(0.4961222873871792, 0.2799572387689701, -0.1270718902659863, 0.3719703269114335, 0.2131214137559209, 0.269285234726165, 2.519088074337557, 9.988238405194851, -8.183023394942618, -0.04335184796346914),
    @(0.4342613726009741, 0.02723224356601772, 0.5262744607434375, 0.3674255475879249, 0.423589368558169, 2.673392208169561, 10.14254253902685, -8.028719261110615, 0.1109522858685349, 2.258522757667927),
    @(-0.6231941776531655, -0.1241519604757457, -0.2830008736312583, -0.2268370526610142, 2.022965786950378, 9.492116117807671, -8.679145682329798, -0.5394741353506483, 1.608096336448744, -1.831933650074586),
    @(0.09201308814246346, -0.06683582501304909, -0.01067200404280504, 2.239130835568587, 9.708281166425881, -8.462980633711588, -0.3233090867324392, 1.824261385066953, -1.615768601456377, -1.676942779870578),
    @(0.3401933040219072, 0.3963571249921513, 2.646159964603544, 10.11531029546084, -8.055951504676631, 0.08372004230251717, 2.231290514101909, -1.20873947242142, -1.269913650835621, 0.7715422863206332),
    @(-0.1026850921852685, 2.147117747426124, 9.616268078283417, -8.554993721854052, -0.4153221748749026, 1.73224829692449, -1.70778168959884, -1.768955868013041, 0.2725000691432133, -0.3100944678904277),
    @(2.305966660581636, 9.77511699143893, -8.396144808698539, -0.2564732617193901, 1.891097210080002, -1.548932776443328, -1.610106954857529, 0.4313489822987259, -0.1512455547349151, -0.05557051392555735)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
